$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A52").Value = 112104863
$ws.Range("B52").Value = 90785
$ws.Range("E52").Value = 1968
$ws.Range("F52").Value = 'Grantaggsvamp'
$ws.Range("G52").Value = 'Bankera violascens'
$ws.Range("H52").Value = '(Alb. & Schwein. : Fr.) Pouzar'
$ws.Range("P52").Value = 'Hökån (Hökån), Jmt'
$ws.Range("Q52").Value = 446637
$ws.Range("R52").Value = 7032524
$ws.Range("A53").Value = 112111386
$ws.Range("B53").Value = 89080
$ws.Range("E53").Value = 256335
$ws.Range("F53").Value = 'Taggfingersvamp'
$ws.Range("G53").Value = 'Ramaria karstenii'
$ws.Range("H53").Value = '(Sacc. & P.Syd.) Corner'
$ws.Range("P53").Value = 'Renkullmyren (Renkullmyren), Jmt'
$ws.Range("Q53").Value = 446734
$ws.Range("R53").Value = 7032709
$ws.Range("AC53").Value = $null
$ws.Range("A54").Value = 112105381
$ws.Range("B54").Value = 89033
$ws.Range("D54").Value = 'NT'
$ws.Range("E54").Value = 3286
$ws.Range("F54").Value = 'Flattoppad klubbsvamp'
$ws.Range("G54").Value = 'Clavariadelphus truncatus'
$ws.Range("H54").Value = '(Quél.) Donk'
$ws.Range("P54").Value = 'Landverktjärnen (Landverktjärnen), Jmt'
$ws.Range("Q54").Value = 446564
$ws.Range("R54").Value = 7032716
$ws.Range("A55").Value = 112111388
$ws.Range("B55").Value = 89090
$ws.Range("D55").Value = 'VU'
$ws.Range("E55").Value = 5747
$ws.Range("F55").Value = 'Läderdoftande fingersvamp'
$ws.Range("G55").Value = 'Ramaria safraniolens'
$ws.Range("H55").Value = 'Christian'
$ws.Range("P55").Value = 'Renkullmyren (Renkullmyren), Jmt'
$ws.Range("Q55").Value = 446734
$ws.Range("R55").Value = 7032709
$ws.Range("A56").Value = 112111398
$ws.Range("B56").Value = 89100
$ws.Range("E56").Value = 5754
$ws.Range("F56").Value = 'Gultoppig fingersvamp'
$ws.Range("G56").Value = 'Ramaria testaceoflava'
$ws.Range("H56").Value = '(Bres.) Corner'
$ws.Range("Q56").Value = 446740
$ws.Range("R56").Value = 7032705
$ws.Range("A57").Value = 112111378
$ws.Range("B57").Value = 83072
$ws.Range("E57").Value = 5589
$ws.Range("F57").Value = 'Rödbrun klubbdyna'
$ws.Range("G57").Value = 'Trichoderma nybergianum'
$ws.Range("H57").Value = '(T.Ulvinen & H.L.Chamb.) Jaklitsch & Voglmayr'
$ws.Range("Q57").Value = 446760
$ws.Range("R57").Value = 7032715
$ws.Range("A58").Value = 112213255
$ws.Range("B58").Value = 90466
$ws.Range("D58").Value = 'LC'
$ws.Range("E58").Value = 4769
$ws.Range("F58").Value = 'Svavelriska'
$ws.Range("G58").Value = 'Lactarius scrobiculatus'
$ws.Range("H58").Value = '(Scop.:Fr.) Fr.'
$ws.Range("P58").Value = 'Ol-olssvarttjärnen, Jmt'
$ws.Range("Q58").Value = 446605
$ws.Range("R58").Value = 7032710
$ws.Range("Y58").Value = '''2023-09-19'
$ws.Range("AA58").Value = '''2023-09-19'
$ws.Range("AW58").Value = 'Erik Lundmark'
$ws.Range("AX58").Value = 'Erik Lundmark'
$ws.Range("B59").Value = 89090
$ws.Range("A60").Value = 112213235
$ws.Range("B60").Value = 90821
$ws.Range("E60").Value = 5964
$ws.Range("F60").Value = 'Fjällig taggsvamp s.str.'
$ws.Range("G60").Value = 'Sarcodon imbricatus s.str.'
$ws.Range("H60").Value = '(L.:Fr.) P.Karst.'
$ws.Range("Q60").Value = 446568
$ws.Range("R60").Value = 7032711
$ws.Range("A61").Value = 112110532
$ws.Range("B61").Value = 88136
$ws.Range("D61").Value = 'VU'
$ws.Range("E61").Value = 245031
$ws.Range("F61").Value = 'Borgsjömusseron'
$ws.Range("G61").Value = 'Tricholoma borgsjoeënse'
$ws.Range("H61").Value = 'Jacobsson & Muskos'
$ws.Range("P61").Value = 'Svensbergsbäcken (Svensbergsbäcken), Jmt'
$ws.Range("Q61").Value = 446765
$ws.Range("R61").Value = 7032863
$ws.Range("Y61").Value = '''2023-09-15'
$ws.Range("AA61").Value = '''2023-09-15'
$ws.Range("AC61").Value = 'På svag sluttning bland kam-och husmossa, revlummer och ekbräken. I närheten finns granvaxskivling, rosa/besk vaxskivling och äggvaxskivling.'
$ws.Range("AW61").Value = 'Rashid Kadhim'
$ws.Range("AX61").Value = 'Rashid Kadhim'

# Structural cell-existence fixes (matching the exact <c> elements added/removed in the diff)
# Row 58 loses its J/K/N/AF placeholder cells entirely (they become fully absent)
$ws.Range("J58").ClearContents()
$ws.Range("K58").ClearContents()
$ws.Range("N58").ClearContents()
$ws.Range("AF58").ClearContents()

# Row 61 gains new empty J/K/N/AF placeholder cells (present but empty, like the rest of rows 52-58)
$ws.Range("J61").Value = "'"
$ws.Range("K61").Value = "'"
$ws.Range("N61").Value = "'"
$ws.Range("AF61").Value = "'"
